$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- Update the re-queried timestamps on the "data" sheet (column F) ---
$newTimes = @(
    "2021-10-05 14:35:28.776854",
    "2021-10-05 14:35:28.776862",
    "2021-10-05 14:35:28.776866",
    "2021-10-05 14:35:28.776868",
    "2021-10-05 14:35:28.776871",
    "2021-10-05 14:35:28.776874",
    "2021-10-05 14:35:28.776877",
    "2021-10-05 14:35:28.776879",
    "2021-10-05 14:35:28.776882",
    "2021-10-05 14:35:28.776884",
    "2021-10-05 14:35:28.776887",
    "2021-10-05 14:35:28.776890",
    "2021-10-05 14:35:28.776892",
    "2021-10-05 14:35:28.776895",
    "2021-10-05 14:35:28.776897",
    "2021-10-05 14:35:28.776900",
    "2021-10-05 14:35:28.776902",
    "2021-10-05 14:35:28.776905",
    "2021-10-05 14:35:28.776908",
    "2021-10-05 14:35:28.776911",
    "2021-10-05 14:35:28.776913",
    "2021-10-05 14:35:28.776916",
    "2021-10-05 14:35:28.776918",
    "2021-10-05 14:35:28.776921"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = 2 + $i
    $data.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add a new "metadata" sheet after "data" ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row (B1:G1) - bold, bordered, centered like the "data" sheet headers
$headers = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $cell = $meta.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# A2 mirrors the styled index column used on the "data" sheet
$a2 = $meta.Cells.Item(2, 1)
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1
$a2.Value = 0

$meta.Cells.Item(2, 2).Value = "Rasopathy"
$meta.Cells.Item(2, 3).Value = 164

$dCell = $meta.Cells.Item(2, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.90"

$meta.Cells.Item(2, 5).Value = "2020-11-10T21:26:55.598088Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:35:28.773173"
$meta.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/164/?format=json"

$data.Activate()
